# Auto-generated PowerShell Excel COM-interop script
# Applies row-level data refresh for rows 20-26 (football fixtures table)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20: A3Rilxx9 | BULGARIA - PARVA LIGA | Lok. Sofia - Botev Plovdiv
$ws.Range("A20").Value = 'A3Rilxx9'
$ws.Range("B20").Value = '18/04/2025'
$ws.Range("C20").Value = '08:00'
$ws.Range("D20").Value = 'BULGARIA - PARVA LIGA'
$ws.Range("E20").Value = 'Lok. Sofia'
$ws.Range("F20").Value = 'Botev Plovdiv'
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 3.6
$ws.Range("I20").Value = 2.2
$ws.Range("J20").Value = 1.05
$ws.Range("K20").Value = 11
$ws.Range("L20").Value = 1.29
$ws.Range("M20").Value = 3.5
$ws.Range("N20").Value = 1.95
$ws.Range("O20").Value = 1.9
$ws.Range("P20").Value = 1.36
$ws.Range("Q20").Value = 3
$ws.Range("R20").Value = 1.67
$ws.Range("S20").Value = 2.1
$ws.Range("T20").Value = 10
$ws.Range("U20").Value = 17
$ws.Range("V20").Value = 11
$ws.Range("W20").Value = 34
$ws.Range("X20").Value = 23
$ws.Range("Y20").Value = 29
$ws.Range("Z20").Value = 12
$ws.Range("AA20").Value = 7
$ws.Range("AB20").Value = 13
$ws.Range("AC20").Value = 41
$ws.Range("AD20").Value = 201
$ws.Range("AE20").Value = 8.5
$ws.Range("AF20").Value = 11
$ws.Range("AG20").Value = 9
$ws.Range("AH20").Value = 21
$ws.Range("AI20").Value = 17
$ws.Range("AJ20").Value = 26

# Row 21: G46WtfUr | BULGARIA - PARVA LIGA | Cherno More - Beroe
$ws.Range("A21").Value = 'G46WtfUr'
$ws.Range("B21").Value = '18/04/2025'
$ws.Range("C21").Value = '10:30'
$ws.Range("D21").Value = 'BULGARIA - PARVA LIGA'
$ws.Range("E21").Value = 'Cherno More'
$ws.Range("F21").Value = 'Beroe'
$ws.Range("G21").Value = 1.83
$ws.Range("H21").Value = 3.2
$ws.Range("I21").Value = 4.75
$ws.Range("J21").Value = 1.1
$ws.Range("K21").Value = 7
$ws.Range("L21").Value = 1.5
$ws.Range("M21").Value = 2.5
$ws.Range("N21").Value = 2.5
$ws.Range("O21").Value = 1.5
$ws.Range("P21").Value = 1.57
$ws.Range("Q21").Value = 2.25
$ws.Range("R21").Value = 2.25
$ws.Range("S21").Value = 1.57
$ws.Range("T21").Value = 5
$ws.Range("U21").Value = 7.5
$ws.Range("V21").Value = 9.5
$ws.Range("W21").Value = 15
$ws.Range("X21").Value = 19
$ws.Range("Y21").Value = 41
$ws.Range("Z21").Value = 6.5
$ws.Range("AA21").Value = 6.5
$ws.Range("AB21").Value = 21
$ws.Range("AC21").Value = 81
$ws.Range("AD21").Value = 101
$ws.Range("AE21").Value = 9.5
$ws.Range("AF21").Value = 21
$ws.Range("AG21").Value = 17
$ws.Range("AH21").Value = 51
$ws.Range("AI21").Value = 41
$ws.Range("AJ21").Value = 51

# Row 22: 2c4QXGiS | BULGARIA - PARVA LIGA | Krumovgrad - CSKA Sofia
$ws.Range("A22").Value = '2c4QXGiS'
$ws.Range("B22").Value = '18/04/2025'
$ws.Range("C22").Value = '13:00'
$ws.Range("D22").Value = 'BULGARIA - PARVA LIGA'
$ws.Range("E22").Value = 'Krumovgrad'
$ws.Range("F22").Value = 'CSKA Sofia'
$ws.Range("G22").Value = 9.5
$ws.Range("H22").Value = 4.75
$ws.Range("I22").Value = 1.33
$ws.Range("J22").Value = 1.06
$ws.Range("K22").Value = 9.5
$ws.Range("L22").Value = 1.3
$ws.Range("M22").Value = 3.4
$ws.Range("N22").Value = 2.03
$ws.Range("O22").Value = 1.83
$ws.Range("P22").Value = 1.4
$ws.Range("Q22").Value = 2.75
$ws.Range("R22").Value = 2.5
$ws.Range("S22").Value = 1.5
$ws.Range("T22").Value = 17
$ws.Range("U22").Value = 41
$ws.Range("V22").Value = 29
$ws.Range("W22").Value = 126
$ws.Range("X22").Value = 81
$ws.Range("Y22").Value = 81
$ws.Range("Z22").Value = 9.5
$ws.Range("AA22").Value = 9.5
$ws.Range("AB22").Value = 29
$ws.Range("AC22").Value = 101
$ws.Range("AD22").Value = 101
$ws.Range("AE22").Value = 5.5
$ws.Range("AF22").Value = 5.5
$ws.Range("AG22").Value = 9.5
$ws.Range("AH22").Value = 8
$ws.Range("AI22").Value = 13
$ws.Range("AJ22").Value = 41

# Row 23: b3dJLtQs | CHILE - LIGA DE PRIMERA | A. Italiano - U. Espanola
$ws.Range("A23").Value = 'b3dJLtQs'
$ws.Range("B23").Value = '18/04/2025'
$ws.Range("C23").Value = '19:00'
$ws.Range("D23").Value = 'CHILE - LIGA DE PRIMERA'
$ws.Range("E23").Value = 'A. Italiano'
$ws.Range("F23").Value = 'U. Espanola'
$ws.Range("G23").Value = 2.35
$ws.Range("H23").Value = 3.4
$ws.Range("I23").Value = 2.9
$ws.Range("J23").Value = 1.04
$ws.Range("K23").Value = 13
$ws.Range("L23").Value = 1.22
$ws.Range("M23").Value = 4
$ws.Range("N23").Value = 1.8
$ws.Range("O23").Value = 2
$ws.Range("P23").Value = 1.36
$ws.Range("Q23").Value = 3
$ws.Range("R23").Value = 1.62
$ws.Range("S23").Value = 2.2
$ws.Range("T23").Value = 9.5
$ws.Range("U23").Value = 12
$ws.Range("V23").Value = 9.5
$ws.Range("W23").Value = 23
$ws.Range("X23").Value = 17
$ws.Range("Y23").Value = 23
$ws.Range("Z23").Value = 12
$ws.Range("AA23").Value = 6.5
$ws.Range("AB23").Value = 12
$ws.Range("AC23").Value = 41
$ws.Range("AD23").Value = 151
$ws.Range("AE23").Value = 11
$ws.Range("AF23").Value = 15
$ws.Range("AG23").Value = 11
$ws.Range("AH23").Value = 29
$ws.Range("AI23").Value = 21
$ws.Range("AJ23").Value = 29

# Row 24: YiuDxRP9 | CHILE - LIGA DE ASCENSO | Deportes Temuco - San Luis
$ws.Range("A24").Value = 'YiuDxRP9'
$ws.Range("B24").Value = '18/04/2025'
$ws.Range("C24").Value = '13:30'
$ws.Range("D24").Value = 'CHILE - LIGA DE ASCENSO'
$ws.Range("E24").Value = 'Deportes Temuco'
$ws.Range("F24").Value = 'San Luis'
$ws.Range("G24").Value = 1.95
$ws.Range("H24").Value = 3.55
$ws.Range("I24").Value = 3.65
$ws.Range("J24").Value = 1.07
$ws.Range("K24").Value = 7.5
$ws.Range("L24").Value = 1.34
$ws.Range("M24").Value = 3.1
$ws.Range("N24").Value = 2
$ws.Range("O24").Value = 1.75
$ws.Range("P24").Value = 1.44
$ws.Range("Q24").Value = 2.65
$ws.Range("R24").Value = 1.87
$ws.Range("S24").Value = 1.85
$ws.Range("T24").Value = 6.7
$ws.Range("U24").Value = 9.25
$ws.Range("V24").Value = 9.25
$ws.Range("W24").Value = 17.5
$ws.Range("X24").Value = 17.5
$ws.Range("Y24").Value = 32
$ws.Range("Z24").Value = 7.5
$ws.Range("AA24").Value = 7.1
$ws.Range("AB24").Value = 17.5
$ws.Range("AC24").Value = 90
$ws.Range("AD24").Value = 900
$ws.Range("AE24").Value = 9.75
$ws.Range("AF24").Value = 20
$ws.Range("AG24").Value = 13.5
$ws.Range("AH24").Value = 55
$ws.Range("AI24").Value = 37
$ws.Range("AJ24").Value = 50

# Row 25: vgKiVrIq | CHILE - LIGA DE ASCENSO | S. Wanderers - Recoleta
$ws.Range("A25").Value = 'vgKiVrIq'
$ws.Range("B25").Value = '18/04/2025'
$ws.Range("C25").Value = '13:30'
$ws.Range("D25").Value = 'CHILE - LIGA DE ASCENSO'
$ws.Range("E25").Value = 'S. Wanderers'
$ws.Range("F25").Value = 'Recoleta'
$ws.Range("G25").Value = 2.07
$ws.Range("H25").Value = 3.2
$ws.Range("I25").Value = 3.65
$ws.Range("J25").Value = 1.09
$ws.Range("K25").Value = 6.6
$ws.Range("L25").Value = 1.4
$ws.Range("M25").Value = 2.8
$ws.Range("N25").Value = 2.2
$ws.Range("O25").Value = 1.62
$ws.Range("P25").Value = 1.45
$ws.Range("Q25").Value = 2.6
$ws.Range("R25").Value = 1.93
$ws.Range("S25").Value = 1.78
$ws.Range("T25").Value = 6.3
$ws.Range("U25").Value = 9.5
$ws.Range("V25").Value = 9.5
$ws.Range("W25").Value = 20
$ws.Range("X25").Value = 19.5
$ws.Range("Y25").Value = 37
$ws.Range("Z25").Value = 6.6
$ws.Range("AA25").Value = 6.5
$ws.Range("AB25").Value = 17.5
$ws.Range("AC25").Value = 100
$ws.Range("AD25").Value = 1000
$ws.Range("AE25").Value = 9
$ws.Range("AF25").Value = 20
$ws.Range("AG25").Value = 13.5
$ws.Range("AH25").Value = 60
$ws.Range("AI25").Value = 40
$ws.Range("AJ25").Value = 55

# Row 26: C6IaTMnd | CHILE - LIGA DE ASCENSO | Curico Unido - CD Santa Cruz
$ws.Range("A26").Value = 'C6IaTMnd'
$ws.Range("B26").Value = '18/04/2025'
$ws.Range("C26").Value = '16:30'
$ws.Range("D26").Value = 'CHILE - LIGA DE ASCENSO'
$ws.Range("E26").Value = 'Curico Unido'
$ws.Range("F26").Value = 'CD Santa Cruz'
$ws.Range("G26").Value = 1.93
$ws.Range("H26").Value = 3.4
$ws.Range("I26").Value = 3.9
$ws.Range("J26").Value = 1.08
$ws.Range("K26").Value = 6.8
$ws.Range("L26").Value = 1.4
$ws.Range("M26").Value = 2.82
$ws.Range("N26").Value = 2.15
$ws.Range("O26").Value = 1.65
$ws.Range("P26").Value = 1.5
$ws.Range("Q26").Value = 2.5
$ws.Range("R26").Value = 1.95
$ws.Range("S26").Value = 1.75
$ws.Range("T26").Value = 6
$ws.Range("U26").Value = 8.75
$ws.Range("V26").Value = 9.25
$ws.Range("W26").Value = 17.5
$ws.Range("X26").Value = 18.5
$ws.Range("Y26").Value = 37
$ws.Range("Z26").Value = 6.8
$ws.Range("AA26").Value = 6.8
$ws.Range("AB26").Value = 18.5
$ws.Range("AC26").Value = 110
$ws.Range("AD26").Value = 1000
$ws.Range("AE26").Value = 9.5
$ws.Range("AF26").Value = 22
$ws.Range("AG26").Value = 14.5
$ws.Range("AH26").Value = 65
$ws.Range("AI26").Value = 45
$ws.Range("AJ26").Value = 55
